# Auto-generated Excel COM-interop script
# Applies the 2026-02-17 22:50 meteocat automatic data/banner refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text-valued cells (timestamps, temperatures, pressures, radiation, etc.) ---
# These values are never numeric-looking on their own (units/letters attached),
# so assigning .Value keeps them as literal text without Excel coercing them to numbers.
$ws.Range('E2').Value = '2026-02-17 22:48:26'
$ws.Range('N2').Value = '-0.1 °C 22:27 TU'
$ws.Range('O2').Value = '2.5 °C'
$ws.Range('E3').Value = '2026-02-17 22:48:29'
$ws.Range('E4').Value = '2026-02-17 22:48:31'
$ws.Range('O4').Value = '9.1 °C'
$ws.Range('E5').Value = '2026-02-17 22:48:34'
$ws.Range('E6').Value = '2026-02-17 22:48:36'
$ws.Range('E7').Value = '2026-02-17 22:48:38'
$ws.Range('O7').Value = '13.9 °C'
$ws.Range('E8').Value = '2026-02-17 22:48:41'
$ws.Range('O8').Value = '10.6 °C'
$ws.Range('E9').Value = '2026-02-17 22:48:43'
$ws.Range('O9').Value = '12.0 °C'
$ws.Range('E10').Value = '2026-02-17 22:48:46'
$ws.Range('E11').Value = '2026-02-17 22:48:48'
$ws.Range('O11').Value = '7.1 °C'
$ws.Range('E12').Value = '2026-02-17 22:48:50'
$ws.Range('N12').Value = '8.6 °C 22:00 TU'
$ws.Range('O12').Value = '12.3 °C'
$ws.Range('E13').Value = '2026-02-17 22:48:53'
$ws.Range('J13').Value = '1018.2 hPa'
$ws.Range('O13').Value = '6.6 °C'
$ws.Range('E14').Value = '2026-02-17 22:48:55'
$ws.Range('K14').Value = '13.4 MJ/m2'
$ws.Range('E15').Value = '2026-02-17 22:48:58'
$ws.Range('O15').Value = '11.7 °C'
$ws.Range('E16').Value = '2026-02-17 22:49:00'
$ws.Range('E17').Value = '2026-02-17 22:49:03'
$ws.Range('E18').Value = '2026-02-17 22:49:05'
$ws.Range('E19').Value = '2026-02-17 22:49:07'
$ws.Range('E20').Value = '2026-02-17 22:49:10'
$ws.Range('O20').Value = '-1.8 °C'
$ws.Range('E21').Value = '2026-02-17 22:49:12'
$ws.Range('J21').Value = '1017.3 hPa'
$ws.Range('O21').Value = '9.4 °C'
$ws.Range('E22').Value = '2026-02-17 22:49:15'
$ws.Range('E23').Value = '2026-02-17 22:49:17'
$ws.Range('M23').Value = '-0.1 °C 22:15 TU'
$ws.Range('E24').Value = '2026-02-17 22:49:19'
$ws.Range('O24').Value = '12.6 °C'
$ws.Range('E25').Value = '2026-02-17 22:49:22'
$ws.Range('E26').Value = '2026-02-17 22:49:24'
$ws.Range('E27').Value = '2026-02-17 22:49:26'
$ws.Range('E28').Value = '2026-02-17 22:49:29'
$ws.Range('E29').Value = '2026-02-17 22:49:31'
$ws.Range('O29').Value = '11.8 °C'
$ws.Range('E30').Value = '2026-02-17 22:49:33'
$ws.Range('J30').Value = '1018.6 hPa'
$ws.Range('N30').Value = '7.5 °C 22:12 TU'
$ws.Range('O30').Value = '11.0 °C'
$ws.Range('E31').Value = '2026-02-17 22:49:35'
$ws.Range('J31').Value = '1018.6 hPa'
$ws.Range('K31').Value = '9.4 MJ/m2'
$ws.Range('E32').Value = '2026-02-17 22:49:38'
$ws.Range('E33').Value = '2026-02-17 22:49:40'
$ws.Range('J33').Value = '1017.7 hPa'
$ws.Range('N33').Value = '2.7 °C 22:23 TU'
$ws.Range('O33').Value = '6.3 °C'
$ws.Range('E34').Value = '2026-02-17 22:49:43'
$ws.Range('E35').Value = '2026-02-17 22:49:45'
$ws.Range('E36').Value = '2026-02-17 22:49:47'
$ws.Range('N36').Value = '8.9 °C 22:19 TU'
$ws.Range('E37').Value = '2026-02-17 22:49:50'
$ws.Range('J37').Value = '1019.4 hPa'
$ws.Range('E38').Value = '2026-02-17 22:49:52'
$ws.Range('E39').Value = '2026-02-17 22:49:54'
$ws.Range('E40').Value = '2026-02-17 22:49:57'
$ws.Range('O40').Value = '9.3 °C'
$ws.Range('E41').Value = '2026-02-17 22:49:59'
$ws.Range('O41').Value = '16.2 °C'
$ws.Range('E42').Value = '2026-02-17 22:50:01'
$ws.Range('E43').Value = '2026-02-17 22:50:04'
$ws.Range('E44').Value = '2026-02-17 22:50:06'
$ws.Range('M44').Value = '0.9 °C 22:27 TU'
$ws.Range('O44').Value = '-2.8 °C'
$ws.Range('E45').Value = '2026-02-17 22:50:09'
$ws.Range('E46').Value = '2026-02-17 22:49:43'
$ws.Range('N46').Value = '9.8 °C 22:29 TU'
$ws.Range('O46').Value = '15.0 °C'

# --- Percentage cells (e.g. "76%") ---
# Assigning a bare "NN%" string directly makes Excel reinterpret it as a number
# formatted as a percentage (changing both stored type and style). To keep these as
# literal text identical to the surrounding cells, stage the text in a scratch column
# forced to Text format, copy only the *value* onto the target cell (so the target
# keeps its own original style/border), then remove the scratch column entirely so the
# sheet dimensions/content are unaffected.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = '76%'
$scratch.Copy() | Out-Null
$ws.Range('H6').PasteSpecial(-4163) | Out-Null
$scratch.Value = '64%'
$scratch.Copy() | Out-Null
$ws.Range('H7').PasteSpecial(-4163) | Out-Null
$scratch.Value = '76%'
$scratch.Copy() | Out-Null
$ws.Range('H10').PasteSpecial(-4163) | Out-Null
$scratch.Value = '64%'
$scratch.Copy() | Out-Null
$ws.Range('H12').PasteSpecial(-4163) | Out-Null
$scratch.Value = '48%'
$scratch.Copy() | Out-Null
$ws.Range('H13').PasteSpecial(-4163) | Out-Null
$scratch.Value = '71%'
$scratch.Copy() | Out-Null
$ws.Range('H14').PasteSpecial(-4163) | Out-Null
$scratch.Value = '64%'
$scratch.Copy() | Out-Null
$ws.Range('H16').PasteSpecial(-4163) | Out-Null
$scratch.Value = '76%'
$scratch.Copy() | Out-Null
$ws.Range('H19').PasteSpecial(-4163) | Out-Null
$scratch.Value = '68%'
$scratch.Copy() | Out-Null
$ws.Range('H20').PasteSpecial(-4163) | Out-Null
$scratch.Value = '41%'
$scratch.Copy() | Out-Null
$ws.Range('H21').PasteSpecial(-4163) | Out-Null
$scratch.Value = '69%'
$scratch.Copy() | Out-Null
$ws.Range('H23').PasteSpecial(-4163) | Out-Null
$scratch.Value = '56%'
$scratch.Copy() | Out-Null
$ws.Range('H27').PasteSpecial(-4163) | Out-Null
$scratch.Value = '68%'
$scratch.Copy() | Out-Null
$ws.Range('H29').PasteSpecial(-4163) | Out-Null
$scratch.Value = '67%'
$scratch.Copy() | Out-Null
$ws.Range('H30').PasteSpecial(-4163) | Out-Null
$scratch.Value = '56%'
$scratch.Copy() | Out-Null
$ws.Range('H34').PasteSpecial(-4163) | Out-Null
$scratch.Value = '62%'
$scratch.Copy() | Out-Null
$ws.Range('H36').PasteSpecial(-4163) | Out-Null
$scratch.Value = '54%'
$scratch.Copy() | Out-Null
$ws.Range('H40').PasteSpecial(-4163) | Out-Null
$scratch.Value = '53%'
$scratch.Copy() | Out-Null
$ws.Range('H41').PasteSpecial(-4163) | Out-Null
$scratch.EntireColumn.Delete() | Out-Null
